{"js": "// Locate the paragraph that starts the \"In every pull request...\" CI/CD blurb,\n// rewrite its text to the new consolidated wording, and remove the following\n// two paragraphs (the old \"static code analysis\" and \"PR accepted/refused\"\n// sentences) whose content has been folded into the rewritten paragraph.\nconst body = context.document.body;\n\nconst results = body.search(\"In every pull request, a github workflow\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target paragraph text\");\n}\n\nconst firstRange = results.items[0];\nconst targetParagraph = firstRange.paragraphs.getFirst();\n\nconst secondParagraph = targetParagraph.getNext();\nconst thirdParagraph = secondParagraph.getNext();\n\nconst newText =\n  \"In every pull request, a github workflow (defined with a yml file in the repo) runs the tests. \" +\n  \"We also use SonarCloud for CI: in every pull request, SonarCloud\\u2019s workflow will perform a static code analysis. \" +\n  \"All the tests and the quality gate must be passed in order for the PR to be accepted (may be accepted/refused automatically or manually).\";\n\n// Replace the whole paragraph's text (keeps the existing run/paragraph formatting).\ntargetParagraph.getRange().insertText(newText, Word.InsertLocation.replace);\n\n// Drop the two paragraphs whose sentences are now merged into targetParagraph.\nsecondParagraph.delete();\nthirdParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Locate the paragraph that starts the \"In every pull request...\" CI/CD blurb,\n# rewrite its text to the new consolidated wording, and remove the following\n# two paragraphs (the old \"static code analysis\" and \"PR accepted/refused\"\n# sentences) whose content has been folded into the rewritten paragraph.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $paras.Item($i).Range.Text\n    if ($text -like \"In every pull request, a github workflow*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find target paragraph (`\"In every pull request, a github workflow...`\")\"\n}\n\n$newText = \"In every pull request, a github workflow (defined with a yml file in the repo) runs the tests. \" + `\n    \"We also use SonarCloud for CI: in every pull request, SonarCloud\" + [char]0x2019 + \"s workflow will perform a static code analysis. \" + `\n    \"All the tests and the quality gate must be passed in order for the PR to be accepted (may be accepted/refused automatically or manually).\"\n\n# Replace the paragraph's text (keeps the paragraph mark / formatting intact).\n$targetParagraph = $d.Paragraphs.Item($targetIndex)\n$targetParagraph.Range.Text = $newText\n\n# The two paragraphs that followed are now redundant - delete them both.\n# After the Text= above, paragraph indices are unchanged, so the paragraph\n# right after our (still single) target paragraph is always at $targetIndex + 1.\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n"}
